$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row from "_old"/"_new" suffixes to the respective input
#    file-version suffixes "_FV2310"/"_FV2404" (the "diff" column is kept
#    as-is).
$newHeaders = @(
    "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310",
    "Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
    "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

# 2. Turn the used range into an Excel Table ("Table1") so the headers are
#    exposed as structured table columns.
$tableRange = $ws.Range("A1:U55")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
$lo.ShowTableStyleRowStripes = $true
$lo.ShowTableStyleColumnStripes = $false
$lo.ShowTableStyleFirstColumn = $false
$lo.ShowTableStyleLastColumn = $false

# 3. Freeze the header row.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
